# Adding some new visuals to see emigration trends
$d = $word.ActiveDocument

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>'

function New-WordPackageXml($bodyInnerXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyInnerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Locate the paragraph that ends with "...compared to the rest." (end of the
# "population density" discussion) - the new content is inserted right after it.
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*compared to the rest.*") {
        $anchorIndex = $i
        break
    }
}
Write-Host "Anchor paragraph index: $anchorIndex"

$anchorPara = $d.Paragraphs($anchorIndex)
$tail = $anchorPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

function Fill-NewParagraph($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $body = '<w:body><w:p><w:pPr>' + $rPr + '</w:pPr>' + $innerXml + '</w:p></w:body>'
    $r.InsertXML((New-WordPackageXml $body))
}

# Paragraph 1: "Examining migration datasets, I have got interesting visuals..."
$p1Runs = '<w:r>' + $rPr + '<w:t>Examining migration datasets, I have got interesting</w:t></w:r>' + `
          '<w:r>' + $rPr + '<w:t xml:space="preserve"> visuals, as a starter I have migration from 1987 to 2023:</w:t></w:r>'
Fill-NewParagraph ($anchorIndex + 1) $p1Runs

# Insert the 2nd new empty paragraph after paragraph 1
$p1 = $d.Paragraphs($anchorIndex + 1)
$t2 = $p1.Range
$t2.Collapse(0)
$t2.InsertParagraphAfter()

# Paragraph 2: "Immigration flows are predominant..." with the (Piola, 2015) citation
$p2Runs = '<w:r>' + $rPr + '<w:t xml:space="preserve">Immigration flows are predominant except in three periods 1987-1991, 1993-1995 and 2010-2014. I have clear that last period was triggered by the Irish economic crisis </w:t></w:r>' + `
          '<w:r>' + $rPr + '<w:t>(</w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r>' + $rPr + '<w:t>Piola</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r>' + $rPr + '<w:t>, 2015)</w:t></w:r>' + `
          '<w:r>' + $rPr + '<w:t>.</w:t></w:r>'
Fill-NewParagraph ($anchorIndex + 2) $p2Runs

# Insert the 3rd new (empty) paragraph after paragraph 2
$p2 = $d.Paragraphs($anchorIndex + 2)
$t3 = $p2.Range
$t3.Collapse(0)
$t3.InsertParagraphAfter()

# The 3rd paragraph should be completely empty (pPr only, no run) just like
# the pre-existing blank paragraphs further down in the document.
Fill-NewParagraph ($anchorIndex + 3) ""

Write-Host "Done. Paragraphs.Count = $($d.Paragraphs.Count)"
